# Weekly update: insert a new record at the top of the "Zapallo italiano"
# date-ordered block (row 549), pushing the existing rows 549-617 down to
# 550-618. This mirrors the author's "Fruta / hortaliza, semanal" refresh,
# which prepends the latest week's price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 549:617 down one row, preserving formatting (row 549 becomes blank)
$ws.Rows.Item(549).Insert()

# Populate the newly inserted row with the latest weekly observation
$ws.Cells.Item(549, 1).Value  = 9
$ws.Cells.Item(549, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(549, 3).Value  = "Metropolitana"
$ws.Cells.Item(549, 4).Value  = 45194
$ws.Cells.Item(549, 5).Value  = 13
$ws.Cells.Item(549, 6).Value  = 100112032
$ws.Cells.Item(549, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(549, 8).Value  = "Sin especificar"
$ws.Cells.Item(549, 9).Value  = "Primera"
$ws.Cells.Item(549, 10).Value = 70
$ws.Cells.Item(549, 11).Value = 14000
$ws.Cells.Item(549, 12).Value = 15000
$ws.Cells.Item(549, 13).Value = 14500
$ws.Cells.Item(549, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(549, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(549, 16).Value = 290
$ws.Cells.Item(549, 17).Value = 50
$ws.Cells.Item(549, 18).Value = "Hortaliza"
